$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.261.09'
$ws.Cells.Item(3, 4).Value = '1.583.42'
$ws.Cells.Item(3, 5).Value = '  -1.06%  '
$ws.Cells.Item(4, 5).Value = '  -0.21%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '209.69'
$ws.Cells.Item(5, 5).Value = '  -0.76%  '
$ws.Cells.Item(6, 5).Value = '  -1.20%  '
$ws.Cells.Item(7, 5).Value = '  -0.19%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.0612'
$ws.Cells.Item(8, 5).Value = '  -1.13%  '
$ws.Cells.Item(9, 5).Value = '  -0.56%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '19.57'
$ws.Cells.Item(11, 5).Value = '  +0.17%  '
$ws.Cells.Item(12, 4).Value = '1.806.99'
$ws.Cells.Item(13, 4).Value = '1.594.87'
$ws.Cells.Item(13, 5).Value = '  -0.54%  '
$ws.Cells.Item(14, 5).Value = '  -0.60%  '
$ws.Cells.Item(15, 5).Value = '  -1.14%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '64.52'
$ws.Cells.Item(16, 5).Value = '  -0.97%  '
$ws.Cells.Item(17, 4).Value = '26.263.80'
$ws.Cells.Item(17, 5).Value = '  -1.62%  '
$ws.Cells.Item(19, 5).Value = '  +0.03%  '
$ws.Cells.Item(20, 5).Value = '  -0.12%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '206.55'
$ws.Cells.Item(21, 5).Value = '  -1.97%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '4.25'
$ws.Cells.Item(22, 5).Value = '  -1.04%  '
$ws.Cells.Item(23, 5).Value = '  -3.10%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '8.83'
$ws.Cells.Item(24, 5).Value = '  -1.57%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '144.57'
$ws.Cells.Item(25, 5).Value = '  +0.56%  '
$ws.Cells.Item(26, 5).Value = '  -0.20%  '
$ws.Cells.Item(27, 5).Value = '  -1.28%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.113'
$ws.Cells.Item(28, 5).Value = '  -0.60%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '15.24'
$ws.Cells.Item(29, 5).Value = '  -0.85%  '
$ws.Cells.Item(31, 5).Value = '  -1.27%  '
$ws.Cells.Item(32, 5).Value = '  -0.87%  '
$ws.Cells.Item(33, 5).Value = '  -0.88%  '
$ws.Cells.Item(34, 5).Value = '  +9.57%  '
$ws.Cells.Item(35, 4).Value = '1.284.08'
$ws.Cells.Item(35, 5).Value = '  -0.88%  '
$ws.Cells.Item(36, 5).Value = '  -0.08%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.604'
$ws.Cells.Item(37, 5).Value = '  -0.52%  '
$ws.Cells.Item(38, 5).Value = '  -1.21%  '
$ws.Cells.Item(39, 5).Value = '  -1.51%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.816'
$ws.Cells.Item(40, 5).Value = '  -0.76%  '
$ws.Cells.Item(41, 5).Value = '  +1.25%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.769'
$ws.Cells.Item(42, 5).Value = '  -1.70%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '2.13'
$ws.Cells.Item(43, 5).Value = '  -2.99%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '62.29'
$ws.Cells.Item(44, 5).Value = '  -1.44%  '
$ws.Cells.Item(45, 4).Value = '1.718.81'
$ws.Cells.Item(45, 5).Value = '  -1.08%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '88.79'
$ws.Cells.Item(46, 5).Value = '  -2.56%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.55'
$ws.Cells.Item(47, 5).Value = '  -0.80%  '
$ws.Cells.Item(48, 5).Value = '  -0.08%  '
$ws.Cells.Item(49, 5).Value = '  -1.67%  '
$ws.Cells.Item(50, 4).Value = '0.0₇0993'
$ws.Cells.Item(50, 5).Value = '  -5.37%  '
$ws.Cells.Item(51, 5).Value = '  +0.00%  '
